# Apply updated horarios (schedule) data for Linea 141 scrape at 09:38:04
# Sheets: LP1912 (sheet1), LP1912-215 (sheet2), 6203-6173 (sheet3)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 09:38:04'
$ws.Range("A3").Value = 'Total filas: 149'
# row 13
$ws.Range("A13").Value = '04:01:01'
$ws.Range("B13").Value = '05:22'
$ws.Range("C13").Value = '23_HERNANDEZ'
$ws.Range("D13").Value = 81
$ws.Range("E13").Value = 'LP1912'
# row 14
$ws.Range("A14").Value = '05:20:00'
$ws.Range("B14").Value = '05:22'
$ws.Range("C14").Value = '14_ABASTO'
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 'LP1912'
# row 116
$ws.Range("A116").Value = '07:46:15'
$ws.Range("B116").Value = '09:33'
$ws.Range("C116").Value = '10_OLMOS'
$ws.Range("D116").Value = 107
$ws.Range("E116").Value = 'LP1912'
# row 117
$ws.Range("A117").Value = '08:50:00'
$ws.Range("B117").Value = '09:33'
$ws.Range("C117").Value = '16_SANTA ANA'
$ws.Range("D117").Value = 43
$ws.Range("E117").Value = 'LP1912'
# row 124
$ws.Range("A124").Value = '09:38:04'
$ws.Range("B124").Value = '09:41'
$ws.Range("C124").Value = '14_ABASTO'
$ws.Range("D124").Value = 3
$ws.Range("E124").Value = 'LP1912'
# row 125
$ws.Range("A125").Value = '08:39:56'
$ws.Range("B125").Value = '09:41'
$ws.Range("C125").Value = '215C_EL PATO'
$ws.Range("D125").Value = 62
$ws.Range("E125").Value = 'LP1912'
# row 126
$ws.Range("A126").Value = '09:38:04'
$ws.Range("B126").Value = '09:41'
$ws.Range("C126").Value = '16_SANTA ANA'
$ws.Range("D126").Value = 3
$ws.Range("E126").Value = 'LP1912'
# row 127
$ws.Range("A127").Value = '07:46:15'
$ws.Range("B127").Value = '09:42'
$ws.Range("C127").Value = '215C_EL PATO'
$ws.Range("D127").Value = 116
$ws.Range("E127").Value = 'LP1912'
# row 128
$ws.Range("A128").Value = '08:21:27'
$ws.Range("B128").Value = '09:43'
$ws.Range("C128").Value = '14_ABASTO'
$ws.Range("D128").Value = 82
$ws.Range("E128").Value = 'LP1912'
# row 129
$ws.Range("A129").Value = '07:59:05'
$ws.Range("B129").Value = '09:44'
$ws.Range("C129").Value = '14_ABASTO'
$ws.Range("D129").Value = 105
$ws.Range("E129").Value = 'LP1912'
# row 130
$ws.Range("A130").Value = '08:57:11'
$ws.Range("B130").Value = '09:45'
$ws.Range("C130").Value = '26_HERNANDEZ'
$ws.Range("D130").Value = 48
$ws.Range("E130").Value = 'LP1912'
# row 131
$ws.Range("A131").Value = '07:59:05'
$ws.Range("B131").Value = '09:52'
$ws.Range("C131").Value = '15_ABASTO'
$ws.Range("D131").Value = 113
$ws.Range("E131").Value = 'LP1912'
# row 132
$ws.Range("A132").Value = '08:57:11'
$ws.Range("B132").Value = '09:53'
$ws.Range("C132").Value = '10_OLMOS'
$ws.Range("D132").Value = 56
$ws.Range("E132").Value = 'LP1912'
# row 133
$ws.Range("A133").Value = '09:38:04'
$ws.Range("B133").Value = '09:59'
$ws.Range("C133").Value = '16_SANTA ANA'
$ws.Range("D133").Value = 21
$ws.Range("E133").Value = 'LP1912'
# row 134
$ws.Range("A134").Value = '09:38:04'
$ws.Range("B134").Value = '10:04'
$ws.Range("C134").Value = '11_ETCHEVERRY'
$ws.Range("D134").Value = 26
$ws.Range("E134").Value = 'LP1912'
# row 135
$ws.Range("A135").Value = '09:38:04'
$ws.Range("B135").Value = '10:05'
$ws.Range("C135").Value = '23_HERNANDEZ'
$ws.Range("D135").Value = 27
$ws.Range("E135").Value = 'LP1912'
# row 136
$ws.Range("A136").Value = '08:21:27'
$ws.Range("B136").Value = '10:12'
$ws.Range("C136").Value = '15_ABASTO'
$ws.Range("D136").Value = 111
$ws.Range("E136").Value = 'LP1912'
# row 137
$ws.Range("A137").Value = '09:38:04'
$ws.Range("B137").Value = '10:13'
$ws.Range("C137").Value = '10_OLMOS'
$ws.Range("D137").Value = 35
$ws.Range("E137").Value = 'LP1912'
# row 138
$ws.Range("A138").Value = '09:38:04'
$ws.Range("B138").Value = '10:21'
$ws.Range("C138").Value = '26_HERNANDEZ'
$ws.Range("D138").Value = 43
$ws.Range("E138").Value = 'LP1912'
# row 139
$ws.Range("A139").Value = '08:39:56'
$ws.Range("B139").Value = '10:22'
$ws.Range("C139").Value = '17_ROMERO'
$ws.Range("D139").Value = 103
$ws.Range("E139").Value = 'LP1912'
# row 140
$ws.Range("A140").Value = '09:38:04'
$ws.Range("B140").Value = '10:24'
$ws.Range("C140").Value = '11_ETCHEVERRY'
$ws.Range("D140").Value = 46
$ws.Range("E140").Value = 'LP1912'
# row 141
$ws.Range("A141").Value = '08:39:56'
$ws.Range("B141").Value = '10:26'
$ws.Range("C141").Value = '215A_EL PATO'
$ws.Range("D141").Value = 107
$ws.Range("E141").Value = 'LP1912'
# row 142
$ws.Range("A142").Value = '09:38:04'
$ws.Range("B142").Value = '10:27'
$ws.Range("C142").Value = '26_HERNANDEZ'
$ws.Range("D142").Value = 49
$ws.Range("E142").Value = 'LP1912'
# row 143
$ws.Range("A143").Value = '08:50:00'
$ws.Range("B143").Value = '10:27'
$ws.Range("C143").Value = '215A_EL PATO'
$ws.Range("D143").Value = 97
$ws.Range("E143").Value = 'LP1912'
# row 144
$ws.Range("A144").Value = '08:50:00'
$ws.Range("B144").Value = '10:42'
$ws.Range("C144").Value = '17_ROMERO'
$ws.Range("D144").Value = 112
$ws.Range("E144").Value = 'LP1912'
# row 145
$ws.Range("A145").Value = '08:50:00'
$ws.Range("B145").Value = '10:44'
$ws.Range("C145").Value = '14_ABASTO'
$ws.Range("D145").Value = 114
$ws.Range("E145").Value = 'LP1912'
# row 146
$ws.Range("A146").Value = '09:38:04'
$ws.Range("B146").Value = '10:52'
$ws.Range("C146").Value = '27_EL RETIRO'
$ws.Range("D146").Value = 74
$ws.Range("E146").Value = 'LP1912'
# row 147
$ws.Range("A147").Value = '09:38:04'
$ws.Range("B147").Value = '11:02'
$ws.Range("C147").Value = '215C_EL PATO'
$ws.Range("D147").Value = 84
$ws.Range("E147").Value = 'LP1912'
# row 148
$ws.Range("A148").Value = '09:38:04'
$ws.Range("B148").Value = '11:07'
$ws.Range("C148").Value = '16_P MOR-167 Y 521'
$ws.Range("D148").Value = 89
$ws.Range("E148").Value = 'LP1912'
# row 149
$ws.Range("A149").Value = '09:38:04'
$ws.Range("B149").Value = '11:12'
$ws.Range("C149").Value = '23_HERNANDEZ'
$ws.Range("D149").Value = 94
$ws.Range("E149").Value = 'LP1912'
# row 150
$ws.Range("A150").Value = '09:38:04'
$ws.Range("B150").Value = '11:20'
$ws.Range("C150").Value = '86_EST CHICA-ESC AGRARIA'
$ws.Range("D150").Value = 102
$ws.Range("E150").Value = 'LP1912'
# row 151
$ws.Range("A151").Value = '09:38:04'
$ws.Range("B151").Value = '11:25'
$ws.Range("C151").Value = '16_P MOR-SANTA ANA'
$ws.Range("D151").Value = 107
$ws.Range("E151").Value = 'LP1912'
# row 152
$ws.Range("A152").Value = '09:38:04'
$ws.Range("B152").Value = '11:27'
$ws.Range("C152").Value = '225_C ROCA-H SUR'
$ws.Range("D152").Value = 109
$ws.Range("E152").Value = 'LP1912'
# row 153
$ws.Range("A153").Value = '09:38:04'
$ws.Range("B153").Value = '11:32'
$ws.Range("C153").Value = '81_EL PELIGRO'
$ws.Range("D153").Value = 114
$ws.Range("E153").Value = 'LP1912'
# row 154
$ws.Range("A154").Value = '09:38:04'
$ws.Range("B154").Value = '11:36'
$ws.Range("C154").Value = '11_ETCHEVERRY'
$ws.Range("D154").Value = 118
$ws.Range("E154").Value = 'LP1912'

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 09:38:04'
$ws.Range("A3").Value = 'Total filas: 20'
# row 25
$ws.Range("A25").Value = '09:38:04'
$ws.Range("B25").Value = '11:02'
$ws.Range("C25").Value = '215C_EL PATO'
$ws.Range("D25").Value = 84
$ws.Range("E25").Value = 'LP1912'

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 09:38:04'
$ws.Range("A3").Value = 'Total filas: 30'
# row 35
$ws.Range("A35").Value = '09:38:04'
$ws.Range("B35").Value = '11:14'
$ws.Range("C35").Value = '215C_LA PLATA'
$ws.Range("D35").Value = 96
$ws.Range("E35").Value = 'L6203'

